# Applies the StoryCards.xlsx update:
#  - Tabelle1 row 21 status changes from "in Arbeit" to "fertig"
#  - Tabelle1 row 22 status changes from "jungfräulich" to "in Arbeit"
#  - Tabelle1 row 21 gets estimated effort (K21=2h), actual effort (L21=3h)
#    and a completion date (M21=2011-10-07), matching the existing date style
#  - The active selection on Tabelle1 moves from D27 to B22

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update status values for rows 21 and 22
$ws.Range("B21").Value = "fertig"
$ws.Range("B22").Value = "in Arbeit"

# Fill in effort + completion date for row 21
$ws.Range("K21").Value = "2h"
$ws.Range("L21").Value = "3h"

# Set the completion date value, then copy the date cell's number format
# from E21 so it reuses the existing date style instead of creating a new one
$ws.Range("M21").Value = 40823
$ws.Range("E21").Copy()
$ws.Range("M21").PasteSpecial(-4122)  # xlPasteFormats

# Move the active selection to B22 on Tabelle1
$ws.Activate()
$ws.Range("B22").Select()
